$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be created in this order so they land at
# index 10 ("TP/A") and index 11 ("ALM") respectively.
$ws.Range("I6").Value = "TP/A"
$ws.Range("H6").Value = "ALM"

# Row 7 additions
$ws.Range("H7").Value = 691
$ws.Range("I7").Formula = "=G7/H7"

# Row 8 additions (full row was previously empty except C8)
$ws.Range("D8").Formula = "=188+35215"
$ws.Range("E8").Formula = "=D8*E4"
$ws.Range("F8").Formula = "=128/E8"
$ws.Range("G8").Formula = "=F8/1000000"
$ws.Range("H8").Value = 1821
$ws.Range("I8").Formula = "=G8/H8"

# Row 9 additions
$ws.Range("H9").Value = 805
$ws.Range("I9").Formula = "=G9/H9"

# Column width adjustments (values chosen so the engine's internal
# pixel-rounding lands on the same stored width as the target OOXML)
$ws.Columns.Item(5).ColumnWidth = 10.92
$ws.Columns.Item(7).ColumnWidth = 16.09
$ws.Columns.Item(9).ColumnWidth = 10.92

# Sheet view adjustments
$null = $ws.Range("E12").Select()
